$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "61.006.22"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  -0.98%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.368.58"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  -1.36%  "
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "406.08"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  -1.36%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "135.84"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  +11.28%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.592"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  +2.48%  "
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  +0.07%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.682"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  +6.25%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.120"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  +5.52%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "42.81"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  +3.42%  "
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  -0.87%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "3.900.61"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  -1.26%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "8.36"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  -0.81%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "19.66"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  +0.65%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.370.91"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "  -0.59%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "61.030.63"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -0.89%  "
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  -0.35%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "11.00"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  +1.09%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0000128"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  +7.41%  "
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  -3.32%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "84.06"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  +9.86%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "310.99"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  +3.80%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.69"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  -2.01%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.13"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  -0.12%  "
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  +11.95%  "
$ws.Cells.Item(27, 2).Value = "Filecoin"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.35"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  +8.79%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "29.49"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  -4.19%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.47"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  -8.00%  "
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  +0.60%  "
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  +0.88%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.999"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  -0.15%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "11.30"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  -1.21%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "41.25"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  -2.87%  "
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  -2.97%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0482"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  +0.09%  "
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  -0.92%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.998"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  -0.03%  "
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  -2.90%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.90"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  -3.99%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.99"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  +0.92%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "136.82"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  +2.40%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.124"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  +0.70%  "
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  +1.70%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.287"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  +1.25%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "16.66"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  -3.44%  "
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  +1.26%  "
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  -2.40%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.121.44"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "  -4.01%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.29"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  -4.10%  "
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  -2.74%  "
